$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow programmatic writes, then re-apply protection at the end.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A38).
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

# Refresh Weight (D) and Percent Change (E) figures for each holding row.
$ws.Range("D2").Value = 0.03647900759310319
$ws.Range("E2").Value = 0.001156069364161905
$ws.Range("D3").Value = 0.02047823755313769
$ws.Range("E3").Value = 0.001168224299065601
$ws.Range("D4").Value = 0.01926940157967585
$ws.Range("E4").Value = 0.005241935483871085
$ws.Range("D5").Value = 0.03774545343387269
$ws.Range("E5").Value = 0.001760563380281743
$ws.Range("D6").Value = 0.03518807628078058
$ws.Range("E6").Value = 0.00156311059007419
$ws.Range("D7").Value = 0.0198746374527217
$ws.Range("E7").Value = 0.001736111111110938
$ws.Range("D8").Value = 0.03657383128774002
$ws.Range("E8").Value = 0.005415944540727802
$ws.Range("D9").Value = 0.02047087657630335
$ws.Range("E9").Value = 0.001708018698310054
$ws.Range("D10").Value = 0.02557939449933626
$ws.Range("E10").Value = -0.00229816147082329
$ws.Range("D11").Value = 0.02327653334387044
$ws.Range("E11").Value = 0.009937410782914213
$ws.Range("D12").Value = 0.05654564385825396
$ws.Range("E12").Value = 0.008206613565049592
$ws.Range("D13").Value = 0.0250089187946748
$ws.Range("E13").Value = 0.005150846210449034
$ws.Range("D14").Value = 0.02728775453963965
$ws.Range("E14").Value = 0.003533026113671234
$ws.Range("D15").Value = 0.03279396968330573
$ws.Range("E15").Value = 0.0003491620111730764
$ws.Range("D16").Value = 0.0191593958703182
$ws.Range("E16").Value = 0.0005175983436851439
$ws.Range("D17").Value = 0.03061195123087034
$ws.Range("E17").Value = -0.008850296401436086
$ws.Range("D18").Value = 0.04244318795523361
$ws.Range("E18").Value = 0.002059496567505814
$ws.Range("D19").Value = 0.1265748592686021
$ws.Range("E19").Value = 0.0006631299734747298
$ws.Range("D20").Value = 0.009165847459806127
$ws.Range("E20").Value = -0.001026167265264277
$ws.Range("D21").Value = 0.01565986480032507
$ws.Range("E21").Value = -0.002125034274746418
$ws.Range("D22").Value = 0.0167196409928915
$ws.Range("E22").Value = 0.007692307692307665
$ws.Range("D23").Value = 0.0163199501743657
$ws.Range("E23").Value = -0.01502219187435982
$ws.Range("D24").Value = 0.02145535611000182
$ws.Range("E24").Value = 0.006044462657514682
$ws.Range("D25").Value = 0.01169602989292692
$ws.Range("E25").Value = 0.01639824304538817
$ws.Range("D26").Value = 0.04185405421901288
$ws.Range("E26").Value = -0.0008244023083262553
$ws.Range("D27").Value = 0.02397628620418502
$ws.Range("E27").Value = -0.0001471093022116232
$ws.Range("D28").Value = 0.04599138326096489
$ws.Range("E28").Value = 0.003314393939393812
$ws.Range("D29").Value = 0.05570623690966563
$ws.Range("E29").Value = 0.009174460431654818
$ws.Range("D30").Value = 0.01274532691708335
$ws.Range("E30").Value = 0.01664447403462055
$ws.Range("D31").Value = 0.02061635810457096
$ws.Range("E31").Value = 0.002305918524212025
$ws.Range("D32").Value = 0.01438329761640835
$ws.Range("E32").Value = 0.005775211017325699
$ws.Range("D33").Value = 0.04184025238744849
$ws.Range("E33").Value = 0.002067183462532096
$ws.Range("D34").Value = 0.01650898414890315
$ws.Range("E34").Value = 0.009861932938856066
$ws.Range("E35").Value = 0.002719469952416009

# Restore sheet protection.
$ws.Protect()
